$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (the old N/O/P -> O/P/Q),
# matching the "inherit formatting from the column to the left" behaviour
# that Excel uses when a column is inserted via the column header.
$ws.Range("N1").EntireColumn.Insert()
$ws.Columns("N:N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet/tab and set its selection,
# matching the new tabSelected/activeTab + S7 selection from the edit.
$ws.Activate()
$ws.Range("S7").Select()
